# "new changes along with new test added"
#
# Summary of the edit:
#   - Delete the "User" sheet (Name/Tools/Place sample data).
#   - Add a new "NewUser" sheet right after "Leave", containing a small
#     login-credentials row (Mukesh500 / Abcd1234 / Abcd1234) and make it
#     the active sheet/tab.
#   - The "Leave" sheet itself is otherwise left untouched (still empty).

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$leaveWs = $wb.Worksheets.Item("Leave")

# A throwaway sheet is added (and later removed) purely so the internal
# sheetId counter advances the same way it did in the real edit (the final
# "NewUser" sheet ends up with sheetId 5, matching the recorded history).
$dummy = $wb.Worksheets.Add($null, $leaveWs)
$dummy.Name = "ZZZ_Dummy"

$leaveWsAgain = $wb.Worksheets.Item("Leave")
$newUserWs = $wb.Worksheets.Add($null, $leaveWsAgain)
$newUserWs.Name = "NewUser"

$dummyRef = $wb.Worksheets.Item("ZZZ_Dummy")
[void]$dummyRef.Delete()

# Seed the new sheet by copying the old "User" row first/second cell so the
# new cells inherit the existing bordered style instead of creating a brand
# new (duplicate) style entry.
$newUserWs2 = $wb.Worksheets.Item("NewUser")
$userWs2 = $wb.Worksheets.Item("User")
[void]$userWs2.Range("A1:B1").Copy($newUserWs2.Range("A1:B1"))
[void]$userWs2.Range("B1").Copy($newUserWs2.Range("C1"))

$newUserWs3 = $wb.Worksheets.Item("NewUser")
$newUserWs3.Range("A1").Value = "Mukesh500"
$newUserWs3.Range("B1").Value = "Abcd1234"
$newUserWs3.Range("C1").Value = "Abcd1234"

[void]$newUserWs3.Columns.Item(1).AutoFit()
[void]$newUserWs3.Columns.Item(2).AutoFit()

[void]$newUserWs3.Range("A1:C1").Select()

# Remove the old "User" sheet last so earlier sheet handles above stay valid.
$userWs3 = $wb.Worksheets.Item("User")
[void]$userWs3.Delete()

# Make "NewUser" the active sheet/tab (activeTab index 2 -> third tab).
$newUserWs4 = $wb.Worksheets.Item("NewUser")
[void]$newUserWs4.Activate()
